# Repro for the commit: update the "İşe Başlangıç Tarihi" (work start date)
# value in Q2 by one month (8/1/2020 -> 9/1/2020), and move the active
# selection from A5 to S2 (the user scrolled right to the phone-number
# column before editing the date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data edit: Q2 date serial 44044 (2020-08-01) -> 44075 (2020-09-01)
$ws.Range("Q2").Value = 44075

# View/selection edit: move the active cell/selection to S2
$ws.Range("S2").Select()
